$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.31"
$ws.Range("E2").Value = "'-2.87%"
$ws.Range("D3").Value = "'31.63"
$ws.Range("E3").Value = "'-2.71%"
$ws.Range("D4").Value = "'5.151"
$ws.Range("E4").Value = "'-3.12%"
$ws.Range("D5").Value = "'0.07472"
$ws.Range("E5").Value = "'0.78%"
$ws.Range("D6").Value = "'7.770"
$ws.Range("E6").Value = "'0.30%"
$ws.Range("D7").Value = "'3.810"
$ws.Range("E7").Value = "'2.40%"
$ws.Range("D8").Value = "'1.648"
$ws.Range("E8").Value = "'4.19%"
$ws.Range("D9").Value = "'0.9236"
$ws.Range("E9").Value = "'0.18%"
$ws.Range("D10").Value = "'0.1709"
$ws.Range("E10").Value = "'2.40%"
$ws.Range("D11").Value = "'0.07528"
$ws.Range("E11").Value = "'1.50%"
$ws.Range("D12").Value = "'0.08002"
$ws.Range("E12").Value = "'0.30%"
$ws.Range("D13").Value = "'0.03016"
$ws.Range("E13").Value = "'-2.89%"
$ws.Range("D14").Value = "'0.09914"
$ws.Range("E14").Value = "'1.03%"
$ws.Range("D15").Value = "'0.001498"
$ws.Range("E15").Value = "'-1.52%"
$ws.Range("D16").Value = "'0.04662"
$ws.Range("E16").Value = "'2.45%"
$ws.Range("D17").Value = "'0.006228"
$ws.Range("E17").Value = "'1.19%"
$ws.Range("E18").Value = "'-0.71%"
$ws.Range("D19").Value = "'2.228"
$ws.Range("E19").Value = "'-0.56%"
$ws.Range("E20").Value = "'0.63%"
$ws.Range("D21").Value = "'0.1349"
$ws.Range("E21").Value = "'2.80%"
$ws.Range("D22").Value = "'4.571"
$ws.Range("E22").Value = "'7.43%"
$ws.Range("D23").Value = "'0.1554"
$ws.Range("E23").Value = "'-5.24%"
$ws.Range("E24").Value = "'-0.17%"
$ws.Range("D25").Value = "'0.004411"
$ws.Range("E25").Value = "'-2.76%"
$ws.Range("D26").Value = "'0.0001402"
$ws.Range("E26").Value = "'19.96%"
$ws.Range("D27").Value = "'0.0001810"
$ws.Range("E27").Value = "'8.65%"
$ws.Range("D39").Value = "'0.01654"
$ws.Range("E39").Value = "'1.91%"
$ws.Range("D40").Value = "'0.04527"
$ws.Range("E40").Value = "'0.65%"
$ws.Range("E41").Value = "'-4.80%"
$ws.Range("D42").Value = "'0.1342"
$ws.Range("E42").Value = "'-1.78%"
$ws.Range("D43").Value = "'0.002063"
$ws.Range("E43").Value = "'-5.27%"
$ws.Range("D44").Value = "'0.01335"
$ws.Range("E44").Value = "'-2.84%"
$ws.Range("D45").Value = "'0.00006090"
$ws.Range("E45").Value = "'1.34%"
$ws.Range("D47").Value = "'0.01227"
$ws.Range("E47").Value = "'-5.60%"
